# Update "想去人数" (F column) values on the 展览, 演出 and 全部类型 sheets
# to reflect newly generated output (gh-pages data refresh).

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 5557
$ws1.Range("F3").Value = 627
$ws1.Range("F4").Value = 12670
$ws1.Range("F5").Value = 308
$ws1.Range("F6").Value = 622
$ws1.Range("F7").Value = 193
$ws1.Range("F8").Value = 377
$ws1.Range("F9").Value = 1179
$ws1.Range("F10").Value = 113

# --- Sheet: 演出 (Performance) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 3

# --- Sheet: 全部类型 (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 5557
$ws4.Range("F3").Value = 627
$ws4.Range("F5").Value = 12670
$ws4.Range("F6").Value = 308
$ws4.Range("F7").Value = 622
$ws4.Range("F8").Value = 193
$ws4.Range("F11").Value = 377
$ws4.Range("F12").Value = 1179
$ws4.Range("F13").Value = 3
$ws4.Range("F14").Value = 113
